$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.411.66"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.848.28"
$ws.Range("E3").Value = "  -0.12%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.90"
$ws.Range("E5").Value = "  -0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6323"
$ws.Range("E6").Value = "  -3.66%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.280.55"
$ws.Range("E8").Value = "  +77.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07582"
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "24.57"
$ws.Range("E11").Value = "  +1.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07712"
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6849"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009969"
$ws.Range("E15").Value = "  +5.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.87"
$ws.Range("E16").Value = "  -1.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.174"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.458.34"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "231.27"
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.50"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9995"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.566"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "155.29"
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1387"
$ws.Range("E25").Value = "  -2.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.438"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.67"
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05811"
$ws.Range("E29").Value = "  -3.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.257"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.130"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.022"
$ws.Range("E32").Value = "  -1.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.457.36"
$ws.Range("E33").Value = "  +72.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.870"
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.156"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7219"
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.593"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.248.60"
$ws.Range("E38").Value = "  +4.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.794"
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01809"
$ws.Range("E40").Value = "  +1.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9012"
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.107"
$ws.Range("E42").Value = "  -2.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9995"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.46"
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.93"
$ws.Range("E45").Value = "  +0.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.318"
$ws.Range("E46").Value = "  -1.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.184"
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4017"
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1125"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05745"
$ws.Range("E51").Value = "  +0.25%  "
